$wb = $excel.ActiveWorkbook

# New "想去人数" (F column) values for rows 2-8, shared by both the
# "展览" and "全部类型" worksheets.
$values = @{
    2 = 1389
    3 = 2231
    4 = 361
    5 = 78
    6 = 6423
    7 = 297
    8 = 121
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $values.Keys) {
        $ws.Range("F$row").Value = $values[$row]
    }
}
